$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows (A2:B21) were reordered (same attribute/type pairs, new order).
$ws.Range("A2").Value = "operation_end_time"
$ws.Range("B2").Value = "datetime"

$ws.Range("A3").Value = "SubProcessID"
$ws.Range("B3").Value = "str"

$ws.Range("A4").Value = "case"
$ws.Range("B4").Value = "str"

$ws.Range("A5").Value = "identifier:id"
$ws.Range("B5").Value = "str"

$ws.Range("A6").Value = "org:resource"
$ws.Range("B6").Value = "str"

$ws.Range("A7").Value = "complete_service_time"
$ws.Range("B7").Value = "str"

$ws.Range("A8").Value = "human_workstation_green_button_pressed"
$ws.Range("B8").Value = "float"

$ws.Range("A9").Value = "parameters"
$ws.Range("B9").Value = "dict"

$ws.Range("A10").Value = "concept:name"
$ws.Range("B10").Value = "str"

$ws.Range("A11").Value = "lifecycle:transition"
$ws.Range("B11").Value = "str"

$ws.Range("A12").Value = "response_status_code"
$ws.Range("B12").Value = "float"

$ws.Range("A13").Value = "case:concept:name"
$ws.Range("B13").Value = "str"

$ws.Range("A14").Value = "time:timestamp"
$ws.Range("B14").Value = "datetime"

$ws.Range("A15").Value = "unsatisfied_condition_description"
$ws.Range("B15").Value = "str"

$ws.Range("A16").Value = "process_model_id"
$ws.Range("B16").Value = "str"

$ws.Range("A17").Value = "lifecycle:state"
$ws.Range("B17").Value = "str"

$ws.Range("A18").Value = "current_task"
$ws.Range("B18").Value = "str"

$ws.Range("A19").Value = "planned_operation_time"
$ws.Range("B19").Value = "str"

$ws.Range("A20").Value = "event_id"
$ws.Range("B20").Value = "str"

$ws.Range("A21").Value = "requested_service_url"
$ws.Range("B21").Value = "str"
